$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6462393333333333
$ws.Range("H2").Value = 1.938718
$ws.Range("I2").Value = 0.03461850536298827
$ws.Range("J2").Value = 0.03461850536298827
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 4.220261333333333
$ws.Range("N2").Value = 12.660784
$ws.Range("O2").Value = 0.6739259863235564
$ws.Range("P2").Value = 0.6739259863235564
$ws.Range("Q2").Value = 2.727298870545777
$ws.Range("R2").Value = 24.545689834912
$ws.Range("S2").Value = 0.0233303103717992
$ws.Range("T2").Value = 0.02333031037179919

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6462393333333333
$ws.Range("H3").Value = 1.938718
$ws.Range("I3").Value = 0.03461850536298827
$ws.Range("J3").Value = 0.03461850536298827
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.296447666666667
$ws.Range("N3").Value = 3.889343
$ws.Range("O3").Value = 0.2070274097896007
$ws.Range("P3").Value = 0.2070274097896007
$ws.Range("Q3").Value = 0.8378154758082221
$ws.Range("R3").Value = 7.540339282274
$ws.Range("S3").Value = 0.007166979496086862
$ws.Range("T3").Value = 0.00716697949608686

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6462393333333333
$ws.Range("H4").Value = 1.938718
$ws.Range("I4").Value = 0.03461850536298827
$ws.Range("J4").Value = 0.03461850536298827
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.745494
$ws.Range("N4").Value = 2.236482
$ws.Range("O4").Value = 0.119046603886843
$ws.Range("P4").Value = 0.119046603886843
$ws.Range("Q4").Value = 0.481767545564
$ws.Range("R4").Value = 4.335907910076
$ws.Range("S4").Value = 0.004121215495102216
$ws.Range("T4").Value = 0.004121215495102215

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.322826
$ws.Range("H5").Value = 45.968478
$ws.Range("I5").Value = 0.8208310864042159
$ws.Range("J5").Value = 0.8208310864042158
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 4.220261333333333
$ws.Range("N5").Value = 12.660784
$ws.Range("O5").Value = 0.6739259863235564
$ws.Range("P5").Value = 0.6739259863235564
$ws.Range("Q5").Value = 64.66633008519467
$ws.Range("R5").Value = 581.996970766752
$ws.Range("S5").Value = 0.5531793995099975
$ws.Range("T5").Value = 0.5531793995099974

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.322826
$ws.Range("H6").Value = 45.968478
$ws.Range("I6").Value = 0.8208310864042159
$ws.Range("J6").Value = 0.8208310864042158
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.296447666666667
$ws.Range("N6").Value = 3.889343
$ws.Range("O6").Value = 0.2070274097896007
$ws.Range("P6").Value = 0.2070274097896007
$ws.Range("Q6").Value = 19.86524201443934
$ws.Range("R6").Value = 178.787178129954
$ws.Range("S6").Value = 0.1699345336930487
$ws.Range("T6").Value = 0.1699345336930487

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.322826
$ws.Range("H7").Value = 45.968478
$ws.Range("I7").Value = 0.8208310864042159
$ws.Range("J7").Value = 0.8208310864042158
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.745494
$ws.Range("N7").Value = 2.236482
$ws.Range("O7").Value = 0.119046603886843
$ws.Range("P7").Value = 0.119046603886843
$ws.Range("Q7").Value = 11.423074846044
$ws.Range("R7").Value = 102.807673614396
$ws.Range("S7").Value = 0.09771715320116971
$ws.Range("T7").Value = 0.09771715320116969

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.698388
$ws.Range("H8").Value = 8.095164
$ws.Range("I8").Value = 0.1445504082327959
$ws.Range("J8").Value = 0.1445504082327959
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 4.220261333333333
$ws.Range("N8").Value = 12.660784
$ws.Range("O8").Value = 0.6739259863235564
$ws.Range("P8").Value = 0.6739259863235564
$ws.Range("Q8").Value = 11.38790253873067
$ws.Range("R8").Value = 102.491122848576
$ws.Range("S8").Value = 0.09741627644175969
$ws.Range("T8").Value = 0.09741627644175969

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.698388
$ws.Range("H9").Value = 8.095164
$ws.Range("I9").Value = 0.1445504082327959
$ws.Range("J9").Value = 0.1445504082327959
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.296447666666667
$ws.Range("N9").Value = 3.889343
$ws.Range("O9").Value = 0.2070274097896007
$ws.Range("P9").Value = 0.2070274097896007
$ws.Range("Q9").Value = 3.498318826361333
$ws.Range("R9").Value = 31.48486943725201
$ws.Range("S9").Value = 0.0299258966004651
$ws.Range("T9").Value = 0.0299258966004651

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.698388
$ws.Range("H10").Value = 8.095164
$ws.Range("I10").Value = 0.1445504082327959
$ws.Range("J10").Value = 0.1445504082327959
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.745494
$ws.Range("N10").Value = 2.236482
$ws.Range("O10").Value = 0.119046603886843
$ws.Range("P10").Value = 0.119046603886843
$ws.Range("Q10").Value = 2.011632063672
$ws.Range("R10").Value = 18.104688573048
$ws.Range("S10").Value = 0.01720823519057111
$ws.Range("T10").Value = 0.01720823519057111
